# TC23_Verify_UserRegistration.xlsx - "Update in CLick event14"
# Remove the intermediate registration-confirmation steps (rows 13-21:
# ExistingaccNObutton click, SCROLL_DOWN, WAIT, Acctype click, WAIT,
# T&CCHeckbox click, WAIT, RegistrationSubmit click, VERIFY_TEXT_PRESENT)
# so that the final two CLICK steps (MyaccountSection, Logout) move up to
# become rows 13 and 14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows("13:21").Delete()

$ws.Activate()
$ws.Range("C11").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
